$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

# Fill in row 42 with the new activity log entry
$ws.Range("B42").Value = 6977
$ws.Range("C42").Value = 43926
$ws.Range("D42").Value = 0.53611111111111109
$ws.Range("E42").Value = 0.059722222222222225
$ws.Range("G42").Value = "Renamed .vho and .sdo files for LogicUnit. Exported .vho, .sdo, .map.summary and .fit.summary files to Documentation folder of LogicUnit"

# Update the view state (scroll position and selection) to match the edited workbook
$window = $excel.ActiveWindow
$window.ScrollRow = 28
$window.ScrollColumn = 1
$ws.Range("A44:G45").Select()
$ws.Range("G45").Activate()
